# "Actualización de espacios entre focalizadores"
# Replace the single data row (row 2) of the "Hoja1" sheet with the new
# course/teacher record, drop the (now stale) hyperlink on the teacher's
# e-mail cell (H2), repoint the metacourse-URL hyperlink (V2) at the new
# course, and leave the selection on the full data row (A2:AC2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 cell values ------------------------------------------------
$ws.Range("A2").Value = 212
$ws.Range("B2").Value = "En línea"
$ws.Range("C2").Value = "Reutilizable"
$ws.Range("D2").Value = "Formativo"
$ws.Range("E2").Value = "EDUC_7117"
$ws.Range("F2").Value = "Fundamentos de la Educación y Teorías Pedagógicas"
$ws.Range("G2").Value = "Pérez Bravo Digna Dionisia"
$ws.Range("H2").Value = "ddperez@utpl.edu.ec"
$ws.Range("I2").Value = "I354288"
$ws.Range("J2").Value = "Loja"
$ws.Range("K2").Value = "Ciencias Sociales, Educación y Humanidades"
$ws.Range("L2").Value = "Maestría en Investigación en Educación"
$ws.Range("M2").Value = "Yunga Godoy Deisi Cecibel"
$ws.Range("N2").Value = "dcyunga@utpl.edu.ec"
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = "Total 180: ACD_40 APE_24 AA_116"
$ws.Range("R2").Value = "Unidad de Formación Disciplinar Avanzada"
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = "Maestría académica"
$ws.Range("U2").Value = "NO"
$ws.Range("V2").Value = "https://utpl.instructure.com/courses/56683"
$ws.Range("X2").Value = "EDUC_7117_META"
$ws.Range("AA2").Value = "EDUC_7117"
$ws.Range("AB2").Value = 45877
$ws.Range("AC2").Value = 45877

# --- Hyperlinks ---------------------------------------------------------
# H2 no longer carries the teacher's personal-site hyperlink (nor its
# special "Hipervínculo" look) now that the row is a different teacher.
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq '$H$2') {
        $h.Delete()
    }
}
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# V2 keeps its hyperlink, just repointed at the new Canvas course.
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq '$V$2') {
        $h.Address = "https://utpl.instructure.com/courses/56683"
    }
}

# --- Selection / view ----------------------------------------------------
$ws.Range("A2:AC2").Select() | Out-Null
